# Apply the "456a3b4" data refresh to both the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets, which contain duplicate copies of the
# same event listing data.

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Only the sheets that actually contain the event rows (identified by
    # the row-4 name matching the old "龙泉" entry) need updating; the
    # "演出" / "本地生活" sheets only have a header row.
    $nameCell = $ws.Cells.Item(4, 3).Value2
    if ($nameCell -eq "龙泉·ACG动湿游戏博览会") {
        # Row 2: 丽水·动漫游戏展 -- "想去人数" 455 -> 456
        $ws.Cells.Item(2, 6).Value = 456

        # Row 3: 丽水·CCAC动漫游戏嘉年华 -- "想去人数" 17 -> 19
        $ws.Cells.Item(3, 6).Value = 19

        # Row 4: rename event and update "想去人数" 1 -> 26
        $ws.Cells.Item(4, 3).Value = "丽水·龙泉ACG动漫游戏博览会"
        $ws.Cells.Item(4, 6).Value = 26

        # Row 5: 丽水·CCAC动漫七夕（回馈展） -- "想去人数" 5 -> 6
        $ws.Cells.Item(5, 6).Value = 6
    }
}
